$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.859.27'
$ws.Range('E2').Value = '  +4.53%  '

$ws.Range('D3').Value = '2.771.89'
$ws.Range('E3').Value = '  +4.76%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '338.42'
$ws.Range('E5').Value = '  +3.75%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '114.82'
$ws.Range('E6').Value = '  +1.93%  '

$ws.Range('E7').Value = '  +4.06%  '

$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.573'
$ws.Range('E9').Value = '  +4.33%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.60'
$ws.Range('E10').Value = '  +4.83%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0857'
$ws.Range('E11').Value = '  +5.24%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.02'
$ws.Range('E12').Value = '  +0.05%  '

$ws.Range('E13').Value = '  +1.86%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.56'
$ws.Range('E14').Value = '  +0.06%  '

$ws.Range('D15').Value = '3.212.91'
$ws.Range('E15').Value = '  +5.08%  '

$ws.Range('D16').Value = '2.781.90'
$ws.Range('E16').Value = '  +5.30%  '

$ws.Range('D17').Value = '51.756.75'
$ws.Range('E17').Value = '  +4.40%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.873'
$ws.Range('E18').Value = '  +1.58%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.19'
$ws.Range('E19').Value = '  +10.32%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.95'
$ws.Range('E20').Value = '  +4.10%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.14'
$ws.Range('E21').Value = '  -1.49%  '

$ws.Range('D22').Value = '0.0₃0972'
$ws.Range('E22').Value = '  +2.47%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '275.13'
$ws.Range('E23').Value = '  +2.54%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '69.69'
$ws.Range('E24').Value = '  +1.00%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.72'
$ws.Range('E25').Value = '  +6.42%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.59'
$ws.Range('E26').Value = '  +1.88%  '

$ws.Range('E27').Value = '  +0.01%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.13'
$ws.Range('E28').Value = '  -0.08%  '

$ws.Range('E29').Value = '  +0.99%  '

$ws.Range('E30').Value = '  +1.62%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '34.50'
$ws.Range('E31').Value = '  -0.68%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '50.15'
$ws.Range('E32').Value = '  +1.02%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.68'
$ws.Range('E33').Value = '  +3.74%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0818'
$ws.Range('E34').Value = '  -0.40%  '

$ws.Range('E35').Value = '  +0.04%  '

$ws.Range('E36').Value = '  +2.74%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.86'
$ws.Range('E37').Value = '  -1.61%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.91'
$ws.Range('E38').Value = '  -0.42%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.19'
$ws.Range('E39').Value = '  +2.26%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0379'
$ws.Range('E40').Value = '  +10.83%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.62'
$ws.Range('E41').Value = '  +23.77%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.35'
$ws.Range('E42').Value = '  +2.22%  '

$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '23.21'
$ws.Range('E43').Value = '  -2.09%  '

$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.115'
$ws.Range('E44').Value = '  +2.95%  '

$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '125.52'
$ws.Range('E45').Value = '  -3.15%  '

$ws.Range('D46').Value = '2.067.79'
$ws.Range('E46').Value = '  +0.32%  '

$ws.Range('E47').Value = '  -0.38%  '

$ws.Range('E48').Value = '  +1.01%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.52'
$ws.Range('E49').Value = '  +5.16%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.86'
$ws.Range('E50').Value = '  -0.62%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '58.95'
$ws.Range('E51').Value = '  +0.32%  '
